$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the salary/wage headers (shorten "Library"/"Librarian" to "Lib") ---
# Column order / meaning is unchanged; only the label text is shortened.
$ws.Range("D1").Value = "LibTraineeWageLow"
$ws.Range("E1").Value = "LibTraineeWageHigh"
$ws.Range("F1").Value = "LibTraineeSalLow"
$ws.Range("G1").Value = "LibTraineeSalHigh"
$ws.Range("H1").Value = "LibISalLow"
$ws.Range("I1").Value = "LibISalHigh"
$ws.Range("J1").Value = "LibIIISalLow"
$ws.Range("K1").Value = "LibIIISalHigh"
# L1 ("WeeklyHours") is unchanged.

# --- Remove stale/incorrect "high" values for LibIIISalHigh (column K) on a handful of rows ---
$ws.Range("K4").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("K26").ClearContents()
$ws.Range("K28").ClearContents()
$ws.Range("K32").ClearContents()
$ws.Range("K42").ClearContents()
$ws.Range("K47").ClearContents()

# --- Unhide the histogram input columns (B, C, D, E) so the backing data/functions are visible ---
$ws.Columns("B:C").Hidden = $false
$ws.Columns("B:C").ColumnWidth = 7.8
$ws.Columns("D").Hidden = $false
$ws.Columns("E").Hidden = $false
$ws.Columns("E").ColumnWidth = 7.8

# --- Update the view: scroll position and active selection ---
$ws.Activate()
$ws.Range("K66").Select()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 5
